# Generate Report for Handoff
#
# Updates the localization-status report:
#  - On the "zh-cn" sheet, rows 7, 9, 10, 11, 13, 14:
#      Priority (col E) set to "ht"
#      Latest Handoff Datetime (col H) updated from 2016-09-03 14:23:44 -> 2016-09-03 14:23:58
#  - On the "de-de" sheet, rows 7, 9, 10, 11, 13, 14:
#      Priority (col E) set to "ht"
#      Latest Handoff Datetime (col H) updated from 2016-09-03 14:23:49 -> 2016-09-03 14:24:05

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 13, 14)

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-09-03 14:23:58"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-09-03 14:24:05"
}
